# Update the "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
# These two sheets carry duplicate event rows; the same updated counts apply
# to each, keyed by row number on that sheet.

$wb = $excel.ActiveWorkbook

# Row(F column) -> new value, for sheet "展览" (sheet1)
$updatesExhibition = @{
    3  = 550
    4  = 1570
    5  = 169
    6  = 135
    7  = 5211
    8  = 184
    9  = 760
    10 = 1057
    12 = 358
    14 = 515
    16 = 6534
    20 = 165
    21 = 68
    22 = 2
    23 = 15613
    24 = 1542
    25 = 10
    26 = 301
    27 = 152
    29 = 11127
    30 = 787
    31 = 4368
    32 = 261
    35 = 310
    36 = 131
}

# Row(F column) -> new value, for sheet "全部类型" (sheet4)
$updatesAllTypes = @{
    3  = 550
    4  = 1570
    5  = 169
    6  = 135
    8  = 5211
    9  = 184
    10 = 760
    12 = 1057
    14 = 358
    16 = 515
    19 = 6534
    23 = 165
    24 = 68
    26 = 2
    27 = 15613
    28 = 1542
    29 = 10
    30 = 301
    31 = 152
    34 = 11127
    35 = 787
    36 = 4368
    37 = 261
    40 = 310
    41 = 131
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibition.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $updatesExhibition[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAllTypes.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $updatesAllTypes[$row]
}
